# financial_model_template.xlsx: correct related data value
# - Cash Flow!B9 / B10 / B18 formula fixes (typo'd reference + exit-IRR range)
# - Active sheet/selection moves from "Inputs" to "Cash Flow"
# - Header / section-label fonts (CJK 宋体, bold) gain explicit Family=3

$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("Inputs")
$wsCashFlow = $wb.Worksheets.Item("Cash Flow")
$wsRentRoll = $wb.Worksheets.Item("Rent Roll")

# --- Formula corrections on the Cash Flow sheet ---------------------------
# Exit NOI (Forward) was pulling next year's rent from the wrong column.
$wsCashFlow.Range("B9").Formula = "=K4*(1+Rent_Growth)"
# Exit Value: same calc, just tidy the spacing to match the fixed sheet.
$wsCashFlow.Range("B10").Formula = "=B9/Exit_Yield"
# Levered IRR needs to include the final (Year 10) cash flow column L.
$wsCashFlow.Range("B18").Formula = "=IRR(B15:L15)"

# --- Font family fix: the bold 宋体 header/label fonts were missing an ---
# --- explicit Family, so set it on every cell that uses those fonts.    ---
$wsInputs.Range("A1:D1").Font.Family = 3
$wsInputs.Range("A12").Font.Family = 3

$wsCashFlow.Range("A1:K1").Font.Family = 3
$wsCashFlow.Range("A4").Font.Family = 3
$wsCashFlow.Range("B4:K4").Font.Family = 3
$wsCashFlow.Range("A6").Font.Family = 3
$wsCashFlow.Range("B6:K6").Font.Family = 3
$wsCashFlow.Range("A8").Font.Family = 3
$wsCashFlow.Range("B12").Font.Family = 3
$wsCashFlow.Range("A14").Font.Family = 3
$wsCashFlow.Range("A17").Font.Family = 3

$wsRentRoll.Range("A1:G1").Font.Family = 3

# --- Active sheet / selection moves from Inputs to Cash Flow --------------
$wsInputs.Activate() | Out-Null
$wsInputs.Range("C15").Select() | Out-Null

$wsCashFlow.Activate() | Out-Null
$wsCashFlow.Range("C19").Select() | Out-Null

Write-Host "edit.ps1 applied"
